$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hand pollinations")

# Insert a new column before column B (shifts Tray#, Page, etc. one column right)
$ws.Columns("B:B").Insert()

# Match the formatting of the column that got pushed to C (old column B),
# which in turn had the same width/style as the former column B.
$ws.Columns("B:B").ColumnWidth = 9

# New header cell
$ws.Range("B1").Value = "Set_num"

# Update the print area to account for the new column
$ws.PageSetup.PrintArea = "`$B`$1:`$N`$100"

# Update the selected cell
$ws.Range("C4").Select()

$wb.Save()
